$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column (D) as Text so that numeric-looking
# values (e.g. "579.66", "9.00") are written verbatim and keep
# trailing zeros / exact formatting instead of being auto-coerced
# into floating-point numbers by the input parser.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '62.144.79'
$ws.Range('E2').Value = '  -3.33%  '

# Row 3
$ws.Range('D3').Value = '2.987.81'
$ws.Range('E3').Value = '  -4.22%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').Value = '579.66'
$ws.Range('E5').Value = '  -2.41%  '

# Row 6
$ws.Range('D6').Value = '145.66'
$ws.Range('E6').Value = '  -7.58%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('E8').Value = '  -3.60%  '

# Row 9
$ws.Range('D9').Value = '2.990.49'
$ws.Range('E9').Value = '  -4.17%  '

# Row 10
$ws.Range('E10').Value = '  -7.01%  '

# Row 11
$ws.Range('E11').Value = '  -5.13%  '

# Row 12
$ws.Range('D12').Value = '0.441'
$ws.Range('E12').Value = '  -2.87%  '

# Row 13
$ws.Range('D13').Value = '0.0000227'
$ws.Range('E13').Value = '  -5.68%  '

# Row 14
$ws.Range('D14').Value = '34.41'
$ws.Range('E14').Value = '  -7.65%  '

# Row 15
$ws.Range('E15').Value = '  +1.35%  '

# Row 16
$ws.Range('D16').Value = '3.482.31'
$ws.Range('E16').Value = '  -4.11%  '

# Row 17
$ws.Range('D17').Value = '7.06'
$ws.Range('E17').Value = '  -2.61%  '

# Row 18
$ws.Range('D18').Value = '62.191.59'
$ws.Range('E18').Value = '  -3.08%  '

# Row 19
$ws.Range('D19').Value = '2.993.08'
$ws.Range('E19').Value = '  -3.92%  '

# Row 20
$ws.Range('D20').Value = '453.69'
$ws.Range('E20').Value = '  -5.84%  '

# Row 21
$ws.Range('D21').Value = '13.82'
$ws.Range('E21').Value = '  -4.84%  '

# Row 22
$ws.Range('D22').Value = '0.676'
$ws.Range('E22').Value = '  -5.40%  '

# Row 23
$ws.Range('D23').Value = '7.27'
$ws.Range('E23').Value = '  -4.66%  '

# Row 24
$ws.Range('D24').Value = '79.87'
$ws.Range('E24').Value = '  -1.81%  '

# Row 25
$ws.Range('D25').Value = '2.26'
$ws.Range('E25').Value = '  -7.80%  '

# Row 26
$ws.Range('D26').Value = '12.24'
$ws.Range('E26').Value = '  -5.72%  '

# Row 27
$ws.Range('E27').Value = '  -0.25%  '

# Row 28
$ws.Range('D28').Value = '9.98'
$ws.Range('E28').Value = '  -4.09%  '

# Row 29
$ws.Range('E29').Value = '  +0.32%  '

# Row 30
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '7.12'
$ws.Range('E30').Value = '  -4.66%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.60'
$ws.Range('E31').Value = '  -3.45%  '

# Row 32
$ws.Range('D32').Value = '2.08'
$ws.Range('E32').Value = '  -5.75%  '

# Row 33
$ws.Range('D33').Value = '26.76'
$ws.Range('E33').Value = '  -2.60%  '

# Row 34
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  -6.27%  '

# Row 35
$ws.Range('E35').Value = '  -4.19%  '

# Row 36
$ws.Range('D36').Value = '0.0₃0781'
$ws.Range('E36').Value = '  -7.34%  '

# Row 37
$ws.Range('D37').Value = '5.71'
$ws.Range('E37').Value = '  -5.58%  '

# Row 38
$ws.Range('D38').Value = '2.10'
$ws.Range('E38').Value = '  -6.17%  '

# Row 39
$ws.Range('D39').Value = '50.10'
$ws.Range('E39').Value = '  -1.95%  '

# Row 40
$ws.Range('D40').Value = '9.00'
$ws.Range('E40').Value = '  -2.39%  '

# Row 41
$ws.Range('E41').Value = '  -11.53%  '

# Row 42
$ws.Range('D42').Value = '407.27'
$ws.Range('E42').Value = '  -8.32%  '

# Row 43
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '0.275'
$ws.Range('E43').Value = '  -6.02%  '

# Row 44
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '0.110'
$ws.Range('E44').Value = '  -1.71%  '

# Row 45
$ws.Range('D45').Value = '2.763.75'
$ws.Range('E45').Value = '  -2.90%  '

# Row 46
$ws.Range('D46').Value = '0.0349'
$ws.Range('E46').Value = '  -4.54%  '

# Row 47
$ws.Range('D47').Value = '37.90'
$ws.Range('E47').Value = '  -5.44%  '

# Row 48
$ws.Range('D48').Value = '127.84'
$ws.Range('E48').Value = '  -1.92%  '

# Row 49
$ws.Range('E49').Value = '  +0.07%  '

# Row 50
$ws.Range('E50').Value = '  -2.46%  '

# Row 51
$ws.Range('D51').Value = '23.64'
$ws.Range('E51').Value = '  -7.65%  '
